$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Text edits inside the "Models" paragraph.
# ---------------------------------------------------------------------------

$old1 = "available models.  Our models were much more complex than we had anticipated.  For example"
$new1 = "available models (www.blenderswap.com).  Our models were much more complex than we had anticipated, so we actually ended up doing two complicated features.  For example"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

$old2 = "300 trees in total and hundreds of vertices per tree.  In addition many of the toys have over 1000 triangle faces.  Though"
$new2 = "300 trees in total and there are over 240,000 triangles among them.   The toys have between 10,000-35,000 triangle faces each.  Though"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the very end of the document up to a
#    new paragraph that replaces the blank paragraph separating the
#    "Models" paragraph from the "Lighting" heading.
# ---------------------------------------------------------------------------

# Locate the blank paragraph right after the "Models" paragraph (the one
# immediately preceding the "Lighting" heading paragraph).
$paraCount = $d.Paragraphs.Count
$modelsParaIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($paraText -eq "Lighting") {
        $modelsParaIndex = $i - 1
        break
    }
}

$blankRange = $d.Paragraphs.Item($modelsParaIndex).Range
$bookmarkParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document><w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$blankRange.InsertXML($bookmarkParaXml) | Out-Null

# ---------------------------------------------------------------------------
# 3) Strip the old "_GoBack" bookmark that used to sit at the very end of
#    the document (end of the "Camera Movement" section), preserving the
#    paragraph's run / proofErr structure exactly.
# ---------------------------------------------------------------------------

$lastParaRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
# Add a temporary trailing paragraph so the target paragraph is no longer
# the very last one in the document (InsertXML behaves differently -
# appending rather than replacing in place - when applied to the final
# paragraph of the document).
$lastParaRange.InsertParagraphAfter() | Out-Null

$targetIndex = $d.Paragraphs.Count - 1
$targetRange = $d.Paragraphs.Item($targetIndex).Range
$cameraParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document><w:body><w:p><w:r><w:t xml:space="preserve">Originally we used only keyboard input to move (i.e. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wasd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for translation and arrow keys for rotation) but we wanted to show the ability to work with mouse events.  We changed from using the arrow keys to using the mouse click and drag.  Although the arrow keys are more intuitive we </w:t></w:r><w:r><w:t>kept the mouse control in order to show our game was capable of using such input.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetRange.InsertXML($cameraParaXml) | Out-Null

# Remove the now-unneeded temporary trailing paragraph.
$dummyRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$dummyRange.Delete() | Out-Null

Write-Output "done"
